$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data, mirroring the layout/style of existing rows (e.g. row 2)
$row = 6

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"

# Column D carries the date style (s="2") used by the other data rows
$ws.Cells.Item($row, 4).Value = 44644
$ws.Range("D$row").Style = $ws.Range("D2").Style
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112052
$ws.Cells.Item($row, 7).Value = "Albahaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 140
$ws.Cells.Item($row, 11).Value = 2500
$ws.Cells.Item($row, 12).Value = 3000
$ws.Cells.Item($row, 13).Value = 2786
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item($row, 16).Value = 464
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
